$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AYKO")

$row = 90

# Columns that look numeric/date-like must be forced to text so they are
# stored the same way as the rest of the sheet (plain text values), not
# auto-converted by Excel into a number or a date serial.
$ws.Cells.Item($row, 1).Value = "'-557"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = "'8/21/2025"
$ws.Cells.Item($row, 2).Style = "Normal"

$ws.Cells.Item($row, 3).Value = "Av Castañares 4621"

$ws.Cells.Item($row, 4).Value = "'8"
$ws.Cells.Item($row, 4).Style = "Normal"

$ws.Cells.Item($row, 5).Value = "Pendiente ADM"
$ws.Cells.Item($row, 6).Value = "AYKO"
$ws.Cells.Item($row, 7).Value = "Pendiente"
$ws.Cells.Item($row, 8).Value = "Colocar columna para pedir traspaso de nodo telecom"

$ws.Cells.Item($row, 9).Value = 1

$ws.Cells.Item($row, 10).Value = "Cambio"
$ws.Cells.Item($row, 11).Value = "Nodo Teco"
$ws.Cells.Item($row, 12).Value = "Pasante"

$ws.Cells.Item($row, 13).Value = -58.470977
$ws.Cells.Item($row, 14).Value = -34.665358

$ws.Cells.Item($row, 15).Value = "Boedo"
$ws.Cells.Item($row, 16).Value = "Capital Sur"
